# Scheduled-runner update: refresh cached market-price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) on a handful of leve
# rows across each crafting-job sheet, per the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 5284.6665
$ws.Range("I53").Value = 250.5
$ws.Range("J53").Value = 7801.75
$ws.Range("K53").Value = 250.5
$ws.Range("L53").Value = 7801.75
$ws.Range("M53").Value = 386.5
$ws.Range("N53").Value = -9075.75
$ws.Range("H106").Value = 11497461
$ws.Range("I106").Value = 47620936
$ws.Range("J106").Value = 3628.0908
$ws.Range("K106").Value = 47620936
$ws.Range("L106").Value = 3628.0908
$ws.Range("M106").Value = -47620305
$ws.Range("N106").Value = -4890.0908
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
$ws.Range("H116").Value = 4093.6667
$ws.Range("J116").Value = 4458.75
$ws.Range("L116").Value = 4458.75
$ws.Range("N116").Value = -11342.75
$ws.Range("H129").Value = 718.9643
$ws.Range("J129").Value = 806.26086
$ws.Range("L129").Value = 2418.78258
$ws.Range("N129").Value = -12418.78258
$ws.Range("H138").Value = 2619.6
$ws.Range("J138").Value = 2619.6
$ws.Range("L138").Value = 7858.799999999999
$ws.Range("N138").Value = -18138.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6482.4136
$ws.Range("I32").Value = 5190.361
$ws.Range("J32").Value = 12112.071
$ws.Range("K32").Value = 5190.361
$ws.Range("L32").Value = 12112.071
$ws.Range("M32").Value = -4903.361
$ws.Range("N32").Value = -12686.071
$ws.Range("H45").Value = 3364.3845
$ws.Range("I45").Value = 2914
$ws.Range("K45").Value = 2914
$ws.Range("M45").Value = -2537
$ws.Range("H97").Value = 41667710
$ws.Range("I97").Value = 850.94446
$ws.Range("K97").Value = 850.94446
$ws.Range("M97").Value = -354.94446
$ws.Range("H102").Value = 1504.5385
$ws.Range("I102").Value = 1449.8182
$ws.Range("K102").Value = 1449.8182
$ws.Range("M102").Value = 172.1818000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 585.4828
$ws.Range("I94").Value = 431.38095
$ws.Range("J94").Value = 990
$ws.Range("K94").Value = 431.38095
$ws.Range("L94").Value = 990
$ws.Range("M94").Value = 19.61905000000002
$ws.Range("N94").Value = -1892
$ws.Range("H134").Value = 4405.2856
$ws.Range("I134").Value = 4127.069
$ws.Range("J134").Value = 5750
$ws.Range("K134").Value = 12381.207
$ws.Range("L134").Value = 17250
$ws.Range("M134").Value = -9846.207000000002
$ws.Range("N134").Value = -22320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3359.6562
$ws.Range("I31").Value = 834.1111
$ws.Range("K31").Value = 834.1111
$ws.Range("M31").Value = -539.1111
$ws.Range("H34").Value = 3359.6562
$ws.Range("I34").Value = 834.1111
$ws.Range("K34").Value = 834.1111
$ws.Range("M34").Value = -632.1111
$ws.Range("H58").Value = 15251.223
$ws.Range("J58").Value = 32559.812
$ws.Range("L58").Value = 32559.812
$ws.Range("N58").Value = -32965.81200000001
$ws.Range("H99").Value = 3001.6553
$ws.Range("I99").Value = 2527.55
$ws.Range("J99").Value = 4055.2222
$ws.Range("K99").Value = 2527.55
$ws.Range("L99").Value = 4055.2222
$ws.Range("M99").Value = -1029.55
$ws.Range("N99").Value = -7051.2222
$ws.Range("H126").Value = 3001.6553
$ws.Range("I126").Value = 2527.55
$ws.Range("J126").Value = 4055.2222
$ws.Range("K126").Value = 7582.650000000001
$ws.Range("L126").Value = 12165.6666
$ws.Range("M126").Value = -5112.650000000001
$ws.Range("N126").Value = -17105.6666
$ws.Range("H134").Value = 1491.3334
$ws.Range("I134").Value = 1306.3636
$ws.Range("K134").Value = 3919.0908
$ws.Range("M134").Value = -1384.0908
$ws.Range("H136").Value = 15251.223
$ws.Range("J136").Value = 32559.812
$ws.Range("L136").Value = 97679.436
$ws.Range("N136").Value = -102779.436
$ws.Range("H137").Value = 40780
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 713.9400000000001
$ws.Range("J131").Value = 741.0217
$ws.Range("L131").Value = 2223.0651
$ws.Range("N131").Value = -12303.0651
$ws.Range("H132").Value = 491
$ws.Range("I132").Value = 491
$ws.Range("K132").Value = 4419
$ws.Range("M132").Value = -1889
$ws.Range("H138").Value = 2813.7273
$ws.Range("I138").Value = 2016.6666
$ws.Range("K138").Value = 6049.9998
$ws.Range("M138").Value = -909.9997999999996
$ws.Range("H141").Value = 3765.8333
$ws.Range("I141").Value = 2010
$ws.Range("K141").Value = 6030
$ws.Range("M141").Value = -850

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 5757.5
$ws.Range("H126").Value = 2771.8865
$ws.Range("I126").Value = 2273.3333
$ws.Range("J126").Value = 3370.15
$ws.Range("K126").Value = 6819.999899999999
$ws.Range("L126").Value = 10110.45
$ws.Range("M126").Value = -4349.999899999999
$ws.Range("N126").Value = -15050.45
$ws.Range("H132").Value = 33257.61
$ws.Range("I132").Value = 6257.231
$ws.Range("J132").Value = 103458.6
$ws.Range("K132").Value = 18771.693
$ws.Range("L132").Value = 310375.8
$ws.Range("M132").Value = -16241.693
$ws.Range("N132").Value = -315435.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 929.4
$ws.Range("I93").Value = 929.4
$ws.Range("K93").Value = 929.4
$ws.Range("M93").Value = 318.6
$ws.Range("H100").Value = 2175
$ws.Range("I100").Value = 1508.8334
$ws.Range("K100").Value = 1508.8334
$ws.Range("M100").Value = -967.8334
$ws.Range("H136").Value = 1580.2069
$ws.Range("I136").Value = 1446.8889
$ws.Range("J136").Value = 3380
$ws.Range("K136").Value = 4340.6667
$ws.Range("L136").Value = 10140
$ws.Range("M136").Value = -1790.6667
$ws.Range("N136").Value = -15240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3924.75
$ws.Range("J15").Value = 3924.75
$ws.Range("L15").Value = 3924.75
$ws.Range("N15").Value = -4500.75
$ws.Range("H54").Value = 20000
$ws.Range("J54").Value = 20000
$ws.Range("L54").Value = 20000
$ws.Range("N54").Value = -21040
$ws.Range("H62").Value = 4001
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4001.5
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4001.5
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5249.5
$ws.Range("H65").Value = 4001
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4001.5
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 20007.5
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -26247.5
$ws.Range("H81").Value = 2567.6
$ws.Range("I81").Value = 1640
$ws.Range("J81").Value = 3495.2
$ws.Range("K81").Value = 3280
$ws.Range("L81").Value = 6990.4
$ws.Range("M81").Value = -2219
$ws.Range("N81").Value = -9112.4
$ws.Range("H84").Value = 2567.6
$ws.Range("I84").Value = 1640
$ws.Range("J84").Value = 3495.2
$ws.Range("K84").Value = 16400
$ws.Range("L84").Value = 34952
$ws.Range("M84").Value = -11096
$ws.Range("N84").Value = -45560
$ws.Range("H107").Value = 55556216
$ws.Range("I107").Value = 100000296
$ws.Range("J107").Value = 1123.125
$ws.Range("K107").Value = 300000888
$ws.Range("L107").Value = 3369.375
$ws.Range("M107").Value = -299998968
$ws.Range("N107").Value = -7209.375
$ws.Range("H113").Value = 1997.2307
$ws.Range("I113").Value = 2080.0833
$ws.Range("J113").Value = 1003
$ws.Range("K113").Value = 6240.249899999999
$ws.Range("L113").Value = 3009
$ws.Range("M113").Value = -4070.249899999999
$ws.Range("N113").Value = -7349
$ws.Range("H132").Value = 1505.2
$ws.Range("I132").Value = 1391.7
$ws.Range("J132").Value = 1656.5333
$ws.Range("K132").Value = 4175.1
$ws.Range("L132").Value = 4969.5999
$ws.Range("M132").Value = -1645.1
$ws.Range("N132").Value = -10029.5999
$ws.Range("H136").Value = 17859186
$ws.Range("I136").Value = 20834148
$ws.Range("J136").Value = 9413.75
$ws.Range("K136").Value = 62502444
$ws.Range("L136").Value = 28241.25
$ws.Range("M136").Value = -62499894
$ws.Range("N136").Value = -33341.25
